$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21: previously blank placeholder row ("Tablet layout" in col A only).
# Becomes the merged "Odd animation..." / "FusedLocationAPI..." entry with real data.
$ws.Range("A21").Value = "Odd animation when map opens, FusedLocationAPI is deprecated"
$ws.Range("B21").Value = 2
$ws.Range("C21").Value = 43199
$ws.Range("D21").Value = "Changed the map so it no longer used deprecated API (changed to FusedLocationProviderClient) and location listener."
$ws.Range("D21").WrapText = $true
$ws.Rows(21).RowHeight = 23

# Row 22: previously blank placeholder row ("Odd animation when map opens" in col A only).
# Becomes the "Property names do not macth" entry with real data.
$ws.Range("A22").Value = "Property names do not macth"
$ws.Range("B22").Value = 1
$ws.Range("C22").Value = 43199
$ws.Range("D22").Value = "Changed names using the conventions (m for member variables, object type start - e.g. btn for button - followed by object description - e.g. btnCall for call button."
$ws.Range("D22").WrapText = $true
$ws.Rows(22).RowHeight = 23

# Row 23: previously blank placeholder row ("Property names do not macth" in col A only).
# Becomes the "Tablet layout" placeholder (still blank data).
$ws.Range("A23").Value = "Tablet layout"

# Row 24: previously blank placeholder row ("FusedLocationAPI is deprecated" in col A only).
# Becomes fully blank (even subject removed).
$ws.Range("A24").ClearContents()

# Update selection/view to match final state
$null = $ws.Range("D24").Select()
